# Update the "想去人数" (interested-count) figures in the 展览 and 全部类型
# sheets to the latest scraped values, as published to gh-pages at 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 2769
    9  = 1473
    13 = 1230
    15 = 380
    19 = 111
    22 = 2706
    23 = 45
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
